$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C6").Value = 116396.2674259275
$ws.Range("D6").Value = 137970.3412021108
$ws.Range("E6").Value = 110386.2103443562
$ws.Range("F6").Value = 142872.0474071934
$ws.Range("G6").Value = 145018.4047542006
$ws.Range("H6").Value = 150935.753186062

$ws.Range("C7").Value = 1063.962674259275
$ws.Range("D7").Value = 1279.703412021108
$ws.Range("E7").Value = 1003.862103443562
$ws.Range("F7").Value = 1328.720474071934
$ws.Range("G7").Value = 1350.184047542006
$ws.Range("H7").Value = 1409.35753186062

$ws.Range("C10").Value = 23129.23919587775
$ws.Range("D10").Value = 18605.9119577439
$ws.Range("E10").Value = 14789.49233281859
$ws.Range("F10").Value = 13964.99431343781
$ws.Range("G10").Value = 15023.9206638157
$ws.Range("H10").Value = 17418.93006482262

$ws.Range("C11").Value = 54.40280804325789
$ws.Range("D11").Value = 49.18499314966178
$ws.Range("E11").Value = 54.64063130993966
$ws.Range("F11").Value = 42.37528297895315
$ws.Range("G11").Value = 39.74040174934412
$ws.Range("H11").Value = 34.32291977657265

$ws.Range("C12").Value = 1050
$ws.Range("D12").Value = 1036
$ws.Range("E12").Value = 1049
$ws.Range("F12").Value = 965
$ws.Range("G12").Value = 961
$ws.Range("H12").Value = 727

$ws.Range("H13").Value = 24

$ws.Range("H14").Value = 24

$ws.Range("E16").Value = 13049.2300994203

$ws.Range("C17").Value = 36.84210526315789
$ws.Range("D17").Value = 50
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 50
$ws.Range("G17").Value = 50
$ws.Range("H17").Value = 54.16666666666666

$ws.Range("C18").Value = 596.4332210050858
$ws.Range("D18").Value = 594.0382383433043
$ws.Range("E18").Value = 595.247704587504
$ws.Range("F18").Value = 594.0382383433043
$ws.Range("G18").Value = 594.0382383433043
$ws.Range("H18").Value = 594.0382383433043

$ws.Range("C19").Value = -12.84622664545977
$ws.Range("D19").Value = -12.84622664545976
$ws.Range("E19").Value = -30.25774552974276
$ws.Range("F19").Value = -12.84622664545977
$ws.Range("G19").Value = -12.84622664545977
$ws.Range("H19").Value = -12.84622664545977

$ws.Range("C20").Value = 109.8537431620315
$ws.Range("D20").Value = 108.3105034806069
$ws.Range("E20").Value = 98.45023933493006
$ws.Range("F20").Value = 108.3105034806069
$ws.Range("G20").Value = 95.17064076873464
$ws.Range("H20").Value = 70.25489790930764

$ws.Range("C21").Value = -6.053924180611452
$ws.Range("D21").Value = -7.359651255399369
$ws.Range("E21").Value = -7.056369499516325
$ws.Range("F21").Value = -7.359651255399371
$ws.Range("G21").Value = -7.050770438412727
$ws.Range("H21").Value = -7.183926726604693

$ws.Range("C22").Value = 157
$ws.Range("D22").Value = 156.285714285706
$ws.Range("E22").Value = 214.6923076923032
$ws.Range("F22").Value = 156.285714285706
$ws.Range("G22").Value = 135.1111111111111
$ws.Range("H22").Value = 126.8461538461458

$ws.Range("C23").Value = 17.91666666666667
$ws.Range("D23").Value = 21
$ws.Range("E23").Value = 13.0769230769213
$ws.Range("F23").Value = 21
$ws.Range("G23").Value = 23.88888888888889
$ws.Range("H23").Value = 25.72727272726852

$ws.Range("C24").Value = 3.103965465147834
$ws.Range("D24").Value = 4.065324406730905
$ws.Range("E24").Value = 4.571365735143691
$ws.Range("F24").Value = 5.11813259573421
$ws.Range("G24").Value = 5.72421252143448
$ws.Range("H24").Value = 6.287523429234587

$ws.Range("C25").Value = 5599.803548733027
$ws.Range("D25").Value = 9140.738657293627
$ws.Range("E25").Value = 3359.114624805226
$ws.Range("F25").Value = 9490.860529085243
$ws.Range("G25").Value = 7501.022486344477
$ws.Range("H25").Value = 5872.323049419249

$ws.Range("C26").Value = 1.224062603548515
$ws.Range("D26").Value = 1.306932940063031
$ws.Range("E26").Value = 1.192194611661512
$ws.Range("F26").Value = 1.412118977783796
$ws.Range("G26").Value = 1.485903071697004
$ws.Range("H26").Value = 1.552677114823654

$ws.Range("C27").Value = 1.160109981149097
$ws.Range("D27").Value = 1.397518435857397
$ws.Range("E27").Value = 1.123678217442283
$ws.Range("F27").Value = 1.649907918806509
$ws.Range("G27").Value = 1.772028747680416
$ws.Range("H27").Value = 2.091466855445514

$ws.Range("C28").Value = 1.242135535128697
$ws.Range("D28").Value = 1.269034596793366
$ws.Range("E28").Value = 1.231584807451516
$ws.Range("F28").Value = 1.3089847460857
$ws.Range("G28").Value = 1.339484412444384
$ws.Range("H28").Value = 1.383830802968821

$ws.Range("C29").Value = 1.941047968957799
$ws.Range("D29").Value = 2.08097565273191
$ws.Range("E29").Value = 1.852069551655872
$ws.Range("F29").Value = 2.302341720339275
$ws.Range("G29").Value = 2.461118401555732
$ws.Range("H29").Value = 2.604060183685443
